# Katalog guncellendi - Cum 14.11.2025 11:50:33,32
# Adds 7 new "Chino Keten Pantolon" (chino linen trousers) product rows
# to Sheet1, one per color, in the "Jeans" category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$fiyat = "360 Tl"
$kategori = "Jeans"
$aciklama = "Slim fit chino keten pantolon, pamuk ve polyester karışımından oluşan orta kalınlığa sahip kumaşı, konforlu bir kullanım vaat eder.Nefes alan kumaşı ile yaz aylarında serin tutar, terletmez.31-38 Beden seçeneği mevcuttur.Ürünümüz serili olarak satılmaktadır."
$stok = "Var"

# Insert 7 fresh rows below the existing data (after row 30) by copying
# the last row's formatting, so the new "aciklama" cells keep the same
# grey/Arial look already used on the existing long-description cells.
for ($i = 0; $i -lt 7; $i++) {
    $ws.Rows("30:30").Copy()
    $ws.Rows("31:31").Insert(-4121)  # xlShiftDown
}
$excel.CutCopyMode = $false

# Row 31 - Antrasit
$ws.Cells.Item(31, 4).Value = "KETENANTRASİT.jpg"
$ws.Cells.Item(31, 2).Value = $fiyat
$ws.Cells.Item(31, 5).Value = $aciklama
$ws.Cells.Item(31, 1).Value = "Chino Keten Pantolon Antrasit"
$ws.Cells.Item(31, 3).Value = $kategori
$ws.Cells.Item(31, 6).Value = $stok

# Row 32 - Bej
$ws.Cells.Item(32, 4).Value = "KETENBEJ.jpg"
$ws.Cells.Item(32, 1).Value = "Chino Keten Pantolon Bej"
$ws.Cells.Item(32, 3).Value = $kategori
$ws.Cells.Item(32, 2).Value = $fiyat
$ws.Cells.Item(32, 5).Value = $aciklama
$ws.Cells.Item(32, 6).Value = $stok

# Row 33 - Gri
$ws.Cells.Item(33, 1).Value = "Chino Keten Pantolon Gri"
$ws.Cells.Item(33, 4).Value = "KETENGRİ.jpg"
$ws.Cells.Item(33, 2).Value = $fiyat
$ws.Cells.Item(33, 3).Value = $kategori
$ws.Cells.Item(33, 5).Value = $aciklama
$ws.Cells.Item(33, 6).Value = $stok

# Row 34 - Mavi
$ws.Cells.Item(34, 1).Value = "Chino Keten Pantolon Mavi"
$ws.Cells.Item(34, 4).Value = "KETENMAVİ.jpg"
$ws.Cells.Item(34, 2).Value = $fiyat
$ws.Cells.Item(34, 3).Value = $kategori
$ws.Cells.Item(34, 5).Value = $aciklama
$ws.Cells.Item(34, 6).Value = $stok

# Row 35 - Siyah
$ws.Cells.Item(35, 1).Value = "Chino Keten Pantolon Siyah"
$ws.Cells.Item(35, 4).Value = "KETENSİYAH.jpg"
$ws.Cells.Item(35, 2).Value = $fiyat
$ws.Cells.Item(35, 3).Value = $kategori
$ws.Cells.Item(35, 5).Value = $aciklama
$ws.Cells.Item(35, 6).Value = $stok

# Row 36 - Taş
$ws.Cells.Item(36, 1).Value = "Chino Keten Pantolon Taş"
$ws.Cells.Item(36, 4).Value = "KETENTAŞ.jpg"
$ws.Cells.Item(36, 2).Value = $fiyat
$ws.Cells.Item(36, 3).Value = $kategori
$ws.Cells.Item(36, 5).Value = $aciklama
$ws.Cells.Item(36, 6).Value = $stok

# Row 37 - Yeşil
$ws.Cells.Item(37, 1).Value = "Chino Keten Pantolon Yeşil"
$ws.Cells.Item(37, 4).Value = "YEŞİLKETEN.jpg"
$ws.Cells.Item(37, 2).Value = $fiyat
$ws.Cells.Item(37, 3).Value = $kategori
$ws.Cells.Item(37, 5).Value = $aciklama
$ws.Cells.Item(37, 6).Value = $stok

# Leave the view scrolled down to the newly added rows, matching where
# the user ended up after typing them in.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("E33").Select() | Out-Null
